# Generate Report for Handoff
#
# A new handoff XLIFF generation pass completed for the
# "a7a2b9a4-d19d-4343-8d8c-c0be3e8440c0.md" file (row 7 on every sheet),
# refreshing its "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
# timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook
$dateFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview = $wb.Worksheets.Item("Overview")
$cellOverview = $wsOverview.Range("G7")
$cellOverview.Value = "2016-09-04 14:45:21"
$cellOverview.NumberFormat = $dateFormat

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$cellZhCn = $wsZhCn.Range("H7")
$cellZhCn.Value = "2016-09-04 14:45:16"
$cellZhCn.NumberFormat = $dateFormat

$wsDeDe = $wb.Worksheets.Item("de-de")
$cellDeDe = $wsDeDe.Range("H7")
$cellDeDe.Value = "2016-09-04 14:45:21"
$cellDeDe.NumberFormat = $dateFormat
